$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update container names (column A) and image names (column B) for rows 3-6
# so that all rows reference the same "nginx:latest" image, with
# sequentially-numbered nginx container names.
$ws.Range("A3").Value = "nginx-container-2"
$ws.Range("B3").Value = "nginx:latest"

$ws.Range("A4").Value = "nginx-container-3"
$ws.Range("B4").Value = "nginx:latest"

$ws.Range("A5").Value = "nginx-container-4"
$ws.Range("B5").Value = "nginx:latest"

$ws.Range("A6").Value = "nginx-container-5"
$ws.Range("B6").Value = "nginx:latest"
